# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 1.455362044514542
    "C2" = 1.655778082260271
    "D2" = 0.7527432677738641
    "E2" = 0.4942365360607697
    "G2" = 4.358119930609447

    "B3" = 3.286832544864788
    "C3" = 1.655778082260271
    "D3" = 0.1494219747398047
    "E3" = 0.4942365360607697
    "G3" = 5.586269137925634

    "B4" = 3.286832544864788
    "C4" = 1.655778082260271
    "D4" = 0.1494219747398047
    "E4" = 0.4942365360607697
    "G4" = 5.586269137925634

    "B5" = 3.286832544864788
    "C5" = 1.655778082260271
    "D5" = 0.7527432677738641
    "E5" = 0.4942365360607697
    "G5" = 6.189590430959694
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
